$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 396
$ws1.Range("F3").Value = 121
$ws1.Range("F4").Value = 1643
$ws1.Range("F5").Value = 19

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 396
$ws4.Range("F3").Value = 121
$ws4.Range("F4").Value = 1643
$ws4.Range("F7").Value = 421
$ws4.Range("F9").Value = 65
$ws4.Range("F10").Value = 0

$wb.Save()
